$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Capabilities": EndPoint column (D) switches from the pCloudy device
# endpoint to the US endpoint, and the remembered selection moves D8 -> D7.
# ---------------------------------------------------------------------------
$capsSheet = $wb.Worksheets.Item("Capabilities")

# D2 / D3 keep their original (quote-prefixed / text-forced) cell style, so
# write the value with a leading apostrophe to preserve that formatting.
$capsSheet.Cells.Item(2, 4).Value = "'https://us.pcloudy.com"
$capsSheet.Cells.Item(3, 4).Value = "'https://us.pcloudy.com"

# Move the remembered selection to D7 on this sheet without permanently
# changing which tab is active in the workbook (the original file keeps
# "DeviceList" as the active/selected tab).
$originalActiveSheet = $wb.ActiveSheet
$capsSheet.Range("D7").Select() | Out-Null
$originalActiveSheet.Activate() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "DeviceList": the NSTRAX10 / NSTRAX15 devices move from Android
# handsets to iOS handsets.
# ---------------------------------------------------------------------------
$deviceSheet = $wb.Worksheets.Item("DeviceList")

# Column B (NSTRAX10 device) first: Device name, then Version (plain /
# quote-prefixed styles respectively).
$deviceSheet.Cells.Item(1, 2).Value = "APPLE_iPhone13mini_iOS_15.2.0_22426"
$deviceSheet.Cells.Item(2, 2).Value = "'15.2.0"

# Column C (NSTRAX15 device) next: Device name, then Version.
$deviceSheet.Cells.Item(1, 3).Value = "APPLE_iPhoneSE2020_iOS_13.6.1_4989f"
$deviceSheet.Cells.Item(2, 3).Value = "'13.6.1"

# Row 3: OperatingSystem (quote-prefixed styles).
$deviceSheet.Cells.Item(3, 2).Value = "'pCloudyIOS"
$deviceSheet.Cells.Item(3, 3).Value = "'pCloudyIOS"
